$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.308.02"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "3.430.85"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +4.44%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.467.98"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "62.300.84"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.572"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.180"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "3.462.65"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.780"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").Value = "2.541.12"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("E51").Value = "  +0.08%  "
